$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 930.1111
$ws.Range("I18").Value = 470.25
$ws.Range("J18").Value = 1298
$ws.Range("K18").Value = 470.25
$ws.Range("L18").Value = 1298
$ws.Range("M18").Value = -186.25
$ws.Range("N18").Value = -1866
$ws.Range("H19").Value = 547.2
$ws.Range("I19").Value = 642.8182
$ws.Range("J19").Value = 430.33334
$ws.Range("K19").Value = 642.8182
$ws.Range("L19").Value = 430.33334
$ws.Range("M19").Value = -467.8182
$ws.Range("N19").Value = -780.33334
$ws.Range("H40").Value = 2095.261
$ws.Range("I40").Value = 2063.9
$ws.Range("J40").Value = 2119.3845
$ws.Range("K40").Value = 2063.9
$ws.Range("L40").Value = 2119.3845
$ws.Range("M40").Value = -1888.9
$ws.Range("N40").Value = -2469.3845
$ws.Range("H70").Value = 1736.1333
$ws.Range("I70").Value = 2211.3333
$ws.Range("J70").Value = 1023.3333
$ws.Range("K70").Value = 6633.999899999999
$ws.Range("L70").Value = 3069.9999
$ws.Range("M70").Value = -6363.999899999999
$ws.Range("N70").Value = -3609.9999
$ws.Range("H73").Value = 1736.1333
$ws.Range("I73").Value = 2211.3333
$ws.Range("J73").Value = 1023.3333
$ws.Range("K73").Value = 6633.999899999999
$ws.Range("L73").Value = 3069.9999
$ws.Range("M73").Value = -5697.999899999999
$ws.Range("N73").Value = -4941.9999
$ws.Range("H76").Value = 3415.4443
$ws.Range("I76").Value = 2998.3333
$ws.Range("J76").Value = 3624
$ws.Range("K76").Value = 2998.3333
$ws.Range("L76").Value = 3624
$ws.Range("M76").Value = -2683.3333
$ws.Range("N76").Value = -4254
$ws.Range("H79").Value = 3415.4443
$ws.Range("I79").Value = 2998.3333
$ws.Range("J79").Value = 3624
$ws.Range("K79").Value = 2998.3333
$ws.Range("L79").Value = 3624
$ws.Range("M79").Value = -1906.3333
$ws.Range("N79").Value = -5808
$ws.Range("H106").Value = 2000
$ws.Range("J106").Value = 2000
$ws.Range("L106").Value = 2000
$ws.Range("N106").Value = -3262
$ws.Range("H113").Value = 2823.762
$ws.Range("I113").Value = 2218.625
$ws.Range("K113").Value = 2218.625
$ws.Range("M113").Value = 1035.375
$ws.Range("H132").Value = 6188.92
$ws.Range("I132").Value = 6590.684
$ws.Range("J132").Value = 4916.6665
$ws.Range("K132").Value = 19772.052
$ws.Range("L132").Value = 14749.9995
$ws.Range("M132").Value = -17242.052
$ws.Range("N132").Value = -19809.9995
$ws.Range("H137").Value = 1347.3529
$ws.Range("I137").Value = 1025.8889
$ws.Range("J137").Value = 1709
$ws.Range("K137").Value = 3077.6667
$ws.Range("L137").Value = 5127
$ws.Range("M137").Value = -527.6666999999998
$ws.Range("N137").Value = -10227
$ws.Range("H138").Value = 2860.9812
$ws.Range("I138").Value = 1859.1
$ws.Range("J138").Value = 3093.9768
$ws.Range("K138").Value = 5577.299999999999
$ws.Range("L138").Value = 9281.930399999999
$ws.Range("M138").Value = -437.2999999999993
$ws.Range("N138").Value = -19561.9304
$ws.Range("H141").Value = 4015.8333
$ws.Range("I141").Value = 2819
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 8457
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = -3277
$ws.Range("N141").Value = -40360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 863.7273
$ws.Range("I2").Value = 857.3333
$ws.Range("J2").Value = 892.5
$ws.Range("K2").Value = 857.3333
$ws.Range("L2").Value = 892.5
$ws.Range("M2").Value = -744.3333
$ws.Range("N2").Value = -1118.5
$ws.Range("H32").Value = 15834.671
$ws.Range("I32").Value = 4444.4907
$ws.Range("J32").Value = 36716.668
$ws.Range("K32").Value = 4444.4907
$ws.Range("L32").Value = 36716.668
$ws.Range("M32").Value = -4157.4907
$ws.Range("N32").Value = -37290.668
$ws.Range("H45").Value = 2399.6667
$ws.Range("I45").Value = 2479.6
$ws.Range("K45").Value = 2479.6
$ws.Range("M45").Value = -2102.6
$ws.Range("H81").Value = 30433.334
$ws.Range("I81").Value = 8500
$ws.Range("K81").Value = 8500
$ws.Range("M81").Value = -7502
$ws.Range("H84").Value = 30433.334
$ws.Range("I84").Value = 8500
$ws.Range("K84").Value = 25500
$ws.Range("M84").Value = -20508
$ws.Range("H116").Value = 863.7273
$ws.Range("I116").Value = 857.3333
$ws.Range("J116").Value = 892.5
$ws.Range("K116").Value = 857.3333
$ws.Range("L116").Value = 892.5
$ws.Range("M116").Value = 1436.6667
$ws.Range("N116").Value = -5480.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 863.7273
$ws.Range("I3").Value = 857.3333
$ws.Range("J3").Value = 892.5
$ws.Range("K3").Value = 857.3333
$ws.Range("L3").Value = 892.5
$ws.Range("M3").Value = -743.3333
$ws.Range("N3").Value = -1120.5
$ws.Range("H86").Value = 1972
$ws.Range("I86").Value = 1855.7142
$ws.Range("J86").Value = 2243.3333
$ws.Range("K86").Value = 1855.7142
$ws.Range("L86").Value = 2243.3333
$ws.Range("M86").Value = -732.7141999999999
$ws.Range("N86").Value = -4489.3333
$ws.Range("H89").Value = 1972
$ws.Range("I89").Value = 1855.7142
$ws.Range("J89").Value = 2243.3333
$ws.Range("K89").Value = 9278.571
$ws.Range("L89").Value = 11216.6665
$ws.Range("M89").Value = -3662.571
$ws.Range("N89").Value = -22448.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 795.44446
$ws.Range("I105").Value = 737
$ws.Range("K105").Value = 737
$ws.Range("M105").Value = 1010
$ws.Range("H132").Value = 2107.5625
$ws.Range("I132").Value = 1657
$ws.Range("J132").Value = 3098.8
$ws.Range("K132").Value = 4971
$ws.Range("L132").Value = 9296.400000000001
$ws.Range("M132").Value = -2441
$ws.Range("N132").Value = -14356.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1905767.8
$ws.Range("J131").Value = 1107.258
$ws.Range("L131").Value = 3321.774
$ws.Range("N131").Value = -13401.774

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 37500
$ws.Range("J68").Value = 37500
$ws.Range("L68").Value = 37500
$ws.Range("N68").Value = -39122
$ws.Range("H71").Value = 37500
$ws.Range("J71").Value = 37500
$ws.Range("L71").Value = 112500
$ws.Range("N71").Value = -120612
$ws.Range("H80").Value = 2729.0278
$ws.Range("J80").Value = 2210.4167
$ws.Range("L80").Value = 2210.4167
$ws.Range("N80").Value = -4206.4167
$ws.Range("H83").Value = 2729.0278
$ws.Range("J83").Value = 2210.4167
$ws.Range("L83").Value = 11052.0835
$ws.Range("N83").Value = -21036.0835
$ws.Range("H113").Value = 1557.625
$ws.Range("I113").Value = 1270.3334
$ws.Range("J113").Value = 1730
$ws.Range("K113").Value = 1270.3334
$ws.Range("L113").Value = 1730
$ws.Range("M113").Value = 899.6666
$ws.Range("N113").Value = -6070
$ws.Range("H122").Value = 1576.6666
$ws.Range("I122").Value = 3001
$ws.Range("J122").Value = 864.5
$ws.Range("K122").Value = 9003
$ws.Range("L122").Value = 2593.5
$ws.Range("M122").Value = -6553
$ws.Range("N122").Value = -7493.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 900.58826
$ws.Range("I22").Value = 470
$ws.Range("J22").Value = 1385
$ws.Range("K22").Value = 470
$ws.Range("L22").Value = 1385
$ws.Range("M22").Value = -175
$ws.Range("N22").Value = -1975
$ws.Range("H27").Value = 900.58826
$ws.Range("I27").Value = 470
$ws.Range("J27").Value = 1385
$ws.Range("K27").Value = 470
$ws.Range("L27").Value = 1385
$ws.Range("M27").Value = -363
$ws.Range("N27").Value = -1599
$ws.Range("H68").Value = 296703.88
$ws.Range("I68").Value = 1250925
$ws.Range("J68").Value = 3097.3845
$ws.Range("K68").Value = 1250925
$ws.Range("L68").Value = 3097.3845
$ws.Range("M68").Value = -1250176
$ws.Range("N68").Value = -4595.3845
$ws.Range("H71").Value = 296703.88
$ws.Range("I71").Value = 1250925
$ws.Range("J71").Value = 3097.3845
$ws.Range("K71").Value = 6254625
$ws.Range("L71").Value = 15486.9225
$ws.Range("M71").Value = -6250881
$ws.Range("N71").Value = -22974.9225

Write-Host "Applied all changes"